# PRACTICA 5: add a new experiment row to the "Tabla2" results table.
# The new row duplicates the previous last row (100,100,10,relu,relu,
# softmax,RMSprop,16,categorical_crossentropy,Earlystop,...,99.12,97.64)
# except for the "Dropout" column, which becomes "(0.2,-)" instead of "-".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row (values + formatting) into the new row right
# below it, then overwrite only the "Dropout" column (N) with the new value.
$ws.Range("D53:P53").Copy($ws.Range("D54:P54"))

$dropout = $ws.Cells.Item(54, 14)
$dropout.ClearFormats()
$dropout.Value = "(0.2,-)"

# Grow the table ("Tabla2") so the new row becomes part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D7:P54"))

# Match the author's final view/selection state after the edit.
$ws.Range("N55").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
